$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.942.60'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '2.792.54'
$ws.Range("E3").Value = '  -1.04%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '358.92'
$ws.Range("E5").Value = '  +0.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.75'
$ws.Range("E6").Value = '  -2.11%  '

$ws.Range("E7").Value = '  -0.58%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -1.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.25'
$ws.Range("E10").Value = '  -1.80%  '

$ws.Range("E11").Value = '  +2.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0852'
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.57'
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.64'
$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").Value = '3.228.63'
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").Value = '2.795.14'
$ws.Range("E16").Value = '  -1.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.944'
$ws.Range("E17").Value = '  +2.14%  '

$ws.Range("D18").Value = '51.879.32'
$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.46'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("E20").Value = '  -1.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.35'
$ws.Range("E21").Value = '  -0.57%  '

$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").Value = '  -1.33%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.19'
$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.22'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.78'
$ws.Range("E25").Value = '  -0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.51'
$ws.Range("E26").Value = '  -2.03%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  +18.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.27'
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.18'
$ws.Range("E30").Value = '  -3.48%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("E31").Value = '  +4.88%  '

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.19'
$ws.Range("E32").Value = '  -0.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.86'
$ws.Range("E33").Value = '  +0.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0465'
$ws.Range("E34").Value = '  -2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0849'
$ws.Range("E35").Value = '  +0.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.18'
$ws.Range("E36").Value = '  -3.92%  '

$ws.Range("E37").Value = '  -0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.85'
$ws.Range("E38").Value = '  +1.93%  '

$ws.Range("E39").Value = '  -2.87%  '

$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("E41").Value = '  +1.31%  '

$ws.Range("E42").Value = '  -1.72%  '

$ws.Range("E43").Value = '  -1.70%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.03'
$ws.Range("E44").Value = '  -5.36%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '119.38'
$ws.Range("E45").Value = '  -3.98%  '

$ws.Range("D46").Value = '2.083.32'
$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("E49").Value = '  -4.65%  '

$ws.Range("E50").Value = '  -1.65%  '

$ws.Range("E51").Value = '  +30.58%  '
